# Update the division-problem answer grid to the new set of values.
# The table is a single 5-column grid; each "answer row" (1, 5, 9, 13, 17)
# holds five populated cells (the three rows beneath each are blank working
# rows), so cells are addressed directly by (row, column) to avoid any
# ambiguity from duplicate old text (e.g. "83÷3=27, 2" appears twice).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "91÷7=13, 0") {
    throw "Unexpected existing text in cell (1,1): $($cell.Range.Text)"
}
$cell.Range.Text = "15÷6=2, 3"

$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "52÷6=8, 4") {
    throw "Unexpected existing text in cell (1,2): $($cell.Range.Text)"
}
$cell.Range.Text = "69÷3=23, 0"

$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "69÷6=11, 3") {
    throw "Unexpected existing text in cell (1,3): $($cell.Range.Text)"
}
$cell.Range.Text = "56÷8=7, 0"

$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "98÷2=49, 0") {
    throw "Unexpected existing text in cell (1,4): $($cell.Range.Text)"
}
$cell.Range.Text = "49÷7=7, 0"

$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "18÷7=2, 4") {
    throw "Unexpected existing text in cell (1,5): $($cell.Range.Text)"
}
$cell.Range.Text = "69÷2=34, 1"

$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "83÷3=27, 2") {
    throw "Unexpected existing text in cell (5,1): $($cell.Range.Text)"
}
$cell.Range.Text = "72÷7=10, 2"

$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "39÷8=4, 7") {
    throw "Unexpected existing text in cell (5,2): $($cell.Range.Text)"
}
$cell.Range.Text = "50÷5=10, 0"

$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "23÷5=4, 3") {
    throw "Unexpected existing text in cell (5,3): $($cell.Range.Text)"
}
$cell.Range.Text = "56÷8=7, 0"

$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "84÷3=28, 0") {
    throw "Unexpected existing text in cell (5,4): $($cell.Range.Text)"
}
$cell.Range.Text = "13÷2=6, 1"

$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "32÷5=6, 2") {
    throw "Unexpected existing text in cell (5,5): $($cell.Range.Text)"
}
$cell.Range.Text = "64÷7=9, 1"

$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "79÷9=8, 7") {
    throw "Unexpected existing text in cell (9,1): $($cell.Range.Text)"
}
$cell.Range.Text = "43÷5=8, 3"

$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "89÷6=14, 5") {
    throw "Unexpected existing text in cell (9,2): $($cell.Range.Text)"
}
$cell.Range.Text = "76÷9=8, 4"

$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "19÷5=3, 4") {
    throw "Unexpected existing text in cell (9,3): $($cell.Range.Text)"
}
$cell.Range.Text = "32÷6=5, 2"

$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "47÷2=23, 1") {
    throw "Unexpected existing text in cell (9,4): $($cell.Range.Text)"
}
$cell.Range.Text = "92÷5=18, 2"

$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "73÷2=36, 1") {
    throw "Unexpected existing text in cell (9,5): $($cell.Range.Text)"
}
$cell.Range.Text = "95÷7=13, 4"

$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "70÷6=11, 4") {
    throw "Unexpected existing text in cell (13,1): $($cell.Range.Text)"
}
$cell.Range.Text = "25÷3=8, 1"

$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "83÷3=27, 2") {
    throw "Unexpected existing text in cell (13,2): $($cell.Range.Text)"
}
$cell.Range.Text = "72÷5=14, 2"

$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "82÷5=16, 2") {
    throw "Unexpected existing text in cell (13,3): $($cell.Range.Text)"
}
$cell.Range.Text = "40÷4=10, 0"

$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "29÷6=4, 5") {
    throw "Unexpected existing text in cell (13,4): $($cell.Range.Text)"
}
$cell.Range.Text = "52÷5=10, 2"

$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "55÷5=11, 0") {
    throw "Unexpected existing text in cell (13,5): $($cell.Range.Text)"
}
$cell.Range.Text = "60÷6=10, 0"

$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "47÷7=6, 5") {
    throw "Unexpected existing text in cell (17,1): $($cell.Range.Text)"
}
$cell.Range.Text = "63÷7=9, 0"

$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "66÷7=9, 3") {
    throw "Unexpected existing text in cell (17,2): $($cell.Range.Text)"
}
$cell.Range.Text = "83÷4=20, 3"

$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "97÷8=12, 1") {
    throw "Unexpected existing text in cell (17,3): $($cell.Range.Text)"
}
$cell.Range.Text = "39÷9=4, 3"

$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "90÷7=12, 6") {
    throw "Unexpected existing text in cell (17,4): $($cell.Range.Text)"
}
$cell.Range.Text = "24÷9=2, 6"

$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "94÷7=13, 3") {
    throw "Unexpected existing text in cell (17,5): $($cell.Range.Text)"
}
$cell.Range.Text = "46÷4=11, 2"

